# Update "想去人数" (want-to-go count) figures on the 展览, 演出 and 全部类型
# sheets to the newly scraped values (gh-pages data regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 101
$wsExhibition.Range("F4").Value = 7340
$wsExhibition.Range("F6").Value = 434
$wsExhibition.Range("F7").Value = 3854
$wsExhibition.Range("F8").Value = 311
$wsExhibition.Range("F9").Value = 547
$wsExhibition.Range("F11").Value = 617
$wsExhibition.Range("F12").Value = 110

# --- 演出 sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 6

# --- 全部类型 sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 101
$wsAll.Range("F4").Value = 6
$wsAll.Range("F5").Value = 7340
$wsAll.Range("F8").Value = 434
$wsAll.Range("F9").Value = 3854
$wsAll.Range("F10").Value = 311
$wsAll.Range("F11").Value = 547
$wsAll.Range("F13").Value = 617
$wsAll.Range("F14").Value = 110
